$d = $word.ActiveDocument

# Locate the paragraph that begins the footer/site-chrome block
# ("Ver no Jupiter Salvar em pdf Salvar em docx"). The paragraph right
# before it (a blank "Normal" paragraph) and the one right after it
# (the "© 2020 ..." copyright notice) belong to the same block and are
# all being removed, since this boilerplate was dropped in the rebuilt
# site export.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    $blankIndex = $targetIndex - 1
    $copyrightIndex = $targetIndex + 1

    $startRange = $d.Paragraphs.Item($blankIndex).Range.Start
    $endRange = $d.Paragraphs.Item($copyrightIndex).Range.End

    $r = $d.Range($startRange, $endRange)
    $r.Delete()
}
